$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: "is very flat" -> "is quite flat"
# ------------------------------------------------------------------
$d.Content.Find.Execute("is very flat", $false, $false, $false, $false, `
    $false, $true, 1, $false, "is quite flat", 2) | Out-Null

# ------------------------------------------------------------------
# Change 2: remove the lone-space run from the otherwise-empty
# paragraph right after "...table is quite flat. </ab>" (leaves the
# paragraph with just its trailing empty run / pilcrow).
# ------------------------------------------------------------------
$blankPara = $d.Paragraphs.Item(14)
Write-Host "blank para text before=[$($blankPara.Range.Text)]"
$blankRange = $blankPara.Range
$spaceRange = $d.Range($blankRange.Start, $blankRange.End - 1)
Write-Host "space range text=[$($spaceRange.Text)]"
$spaceRange.Delete()

# ------------------------------------------------------------------
# Change 3: "; one says that they are excellent for molding." ->
#           ". They are said to be excellent for molding."
# split across 3 runs:
#   ". T" | "hey are said to be " | "excellent for molding."
# ------------------------------------------------------------------

# 3a. drop the "; " that sits right before "one says..."
$rngAnchor = $d.Content
$fAnchor = $rngAnchor.Find
$fAnchor.ClearFormatting()
$fAnchor.Text = "shells</m>; one"
$fAnchor.Execute() | Out-Null
$semiStart = $rngAnchor.Start + ("shells</m>").Length
$semiSpace = $d.Range($semiStart, $semiStart + 2)
$semiSpace.Delete()

# 3b. rewrite "one says that they are" -> ". They are said to be "
#     (this run already carries the plain / automatic-color formatting
#     we need for the first two new runs)
$rngOne = $d.Content
$fOne = $rngOne.Find
$fOne.ClearFormatting()
$fOne.Text = "one says that they are"
$fOne.Execute() | Out-Null
$oneStart = $rngOne.Start
$rngOne.Text = ". They are said to be "

# split that run into ". T" / "hey are said to be " by toggling Bold
# on then off over the first 3 characters (forces a run break without
# leaving any explicit rPr behind)
$splitRng = $d.Range($oneStart, $oneStart + 3)
$splitRng.Font.Bold = $true
$splitRng2 = $d.Range($oneStart, $oneStart + 3)
$splitRng2.Font.Bold = $false

# 3c. trim the leading space off " excellent for molding." so it butts
#     up against the preceding run (that run already has color 000000,
#     exactly matching the target's third run)
$rngExc = $d.Content
$fExc = $rngExc.Find
$fExc.ClearFormatting()
$fExc.Text = " excellent for molding."
$fExc.Execute() | Out-Null
$leadSpace = $d.Range($rngExc.Start, $rngExc.Start + 1)
$leadSpace.Delete()
